# Fix: dependency rows under "root -> Applications -> UPMIS" were showing the
# parent node's own CI_Name/Dependency_Descrip text ("root" / "Unclaimed
# Property Management System") instead of each dependency's own short name
# and description. Also normalize "root" -> "Root" in CI_Type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network_diagram")

# --- CI_Type: "root" -> "Root" for the application-dependency block (rows 2-9)
$ws.Range("A2:A9").Value = "Root"

# --- Dependency_Name (col F): shorten the "(ABBR)" style names to just the
# abbreviation, for the rows that had one ---
$ws.Range("F3").Value = "ALIS"
$ws.Range("F4").Value = "`tEMILI"
$ws.Range("F7").Value = "`tCCAS"
$ws.Range("F8").Value = "`tLSOP"

# --- Dependency_Descrip (col G): replace the incorrect parent-node text
# ("Unclaimed Property Management System") with each dependency's own name ---
$ws.Range("G3").Value = "Automated Licensing Information System"
$ws.Range("G4").Value = "Electronic Management of Investigative & Licensing Information"
$ws.Range("G7").Value = "Coverage and Compliance Automated System"
$ws.Range("G8").Value = "Legal Service of Process"

# --- Row 9: AIMS ---
$ws.Range("F9").Value = "`tAIMS"
$ws.Range("G9").Value = "Automated Investigative Management System "

# --- Rows 5 & 6 (FLAIR / FCDICE): Dependency_Name already correct, only the
# Dependency_Descrip needed fixing ---
$ws.Range("G6").Value = "Florida State Fire College Electronic Information Database "
$ws.Range("G5").Value = "Florida Accounting Information Resource "

# --- Move the saved cursor/selection from the bottom of the sheet back up near
# the top (matches the author re-reviewing the fixed rows before saving) ---
$ws.Range("B5").Select()
